$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final roster table (player, position, team) for rows 2-18.
$data = @(
    @("Donovan Mitchell", "PG,SG", "Cleveland Cavaliers"),
    @("Kelly Oubre Jr.", "SG,SF", "Philadelphia 76ers"),
    @("Ty Jerome", "PG,SG", "Cleveland Cavaliers"),
    @("Kel'el Ware", "C", "Miami Heat"),
    @("Kristaps Porzingis", "PF,C", "Boston Celtics"),
    @("Victor Wembanyama", "C", "San Antonio Spurs"),
    @("Domantas Sabonis", "C", "Sacramento Kings"),
    @("De'Andre Hunter", "SF,PF", "Atlanta Hawks"),
    @("Andrew Wiggins", "SF,PF", "Golden State Warriors"),
    @("Michael Porter Jr.", "SF,PF", "Denver Nuggets"),
    @("Dyson Daniels", "PG,SG,SF", "Atlanta Hawks"),
    @("Malik Beasley", "SG,SF", "Detroit Pistons"),
    @("Jaden McDaniels", "SF,PF", "Minnesota Timberwolves"),
    @("Alperen Sengün", "C", "Houston Rockets"),
    @("Josh Hart", "SG,SF,PF", "New York Knicks"),
    @("Cam Thomas", "SG,SF", "Brooklyn Nets"),
    @("Donte DiVincenzo", "PG,SG,SF", "Minnesota Timberwolves")
)

$row = 2
foreach ($entry in $data) {
    $ws.Range("A$row").Value = $entry[0]
    $ws.Range("B$row").Value = $entry[1]
    $ws.Range("C$row").Value = $entry[2]
    $row++
}
